{"js": "// Update the \"two-digit number divided by one-digit number\" drill sheet:\n// replace each division-problem cell's text with its new value, in order.\n// Replacements are positional (row/col in the single body table) because a\n// couple of source values repeat (e.g. \"43\u00f79=\" appears twice) but map to\n// different targets, so a global text search/replace would be ambiguous.\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Row indices (0-based) that actually contain problems; the intervening\n// rows are blank spacer rows.\nconst dataRows = [0, 4, 8, 12, 16];\n\n// New values, 5 per data row, in reading order (left-to-right, top-to-bottom).\nconst replacements = [\n  [\"98\u00f79=\", \"26\u00f79=\", \"76\u00f76=\", \"89\u00f76=\", \"39\u00f77=\"],\n  [\"61\u00f78=\", \"22\u00f74=\", \"96\u00f77=\", \"34\u00f79=\", \"93\u00f74=\"],\n  [\"72\u00f72=\", \"27\u00f79=\", \"96\u00f73=\", \"55\u00f75=\", \"90\u00f77=\"],\n  [\"17\u00f73=\", \"52\u00f78=\", \"28\u00f75=\", \"13\u00f79=\", \"68\u00f76=\"],\n  [\"92\u00f78=\", \"89\u00f72=\", \"84\u00f73=\", \"52\u00f75=\", \"11\u00f75=\"],\n];\n\nfor (let r = 0; r < dataRows.length; r++) {\n  const row = dataRows[r];\n  for (let col = 0; col < replacements[r].length; col++) {\n    table.getCell(row, col).value = replacements[r][col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the \"two-digit number divided by one-digit number\" drill sheet:\n# replace each division-problem cell's text with its new value, in order.\n# Replacements are positional (row/col in the single body table) because a\n# couple of source values repeat (e.g. \"43\u00f79=\" appears twice) but map to\n# different targets, so a global text search/replace would be ambiguous.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Row indices (1-based) that actually contain problems; the intervening\n# rows are blank spacer rows.\n$dataRows = @(1, 5, 9, 13, 17)\n\n# New values, 5 per data row, in reading order (left-to-right, top-to-bottom).\n$replacements = @(\n    @(\"98\u00f79=\", \"26\u00f79=\", \"76\u00f76=\", \"89\u00f76=\", \"39\u00f77=\"),\n    @(\"61\u00f78=\", \"22\u00f74=\", \"96\u00f77=\", \"34\u00f79=\", \"93\u00f74=\"),\n    @(\"72\u00f72=\", \"27\u00f79=\", \"96\u00f73=\", \"55\u00f75=\", \"90\u00f77=\"),\n    @(\"17\u00f73=\", \"52\u00f78=\", \"28\u00f75=\", \"13\u00f79=\", \"68\u00f76=\"),\n    @(\"92\u00f78=\", \"89\u00f72=\", \"84\u00f73=\", \"52\u00f75=\", \"11\u00f75=\")\n)\n\nfor ($r = 0; $r -lt $dataRows.Length; $r++) {\n    $row = $dataRows[$r]\n    $vals = $replacements[$r]\n    for ($c = 1; $c -le $vals.Length; $c++) {\n        $t.Cell($row, $c).Range.Text = $vals[$c - 1]\n    }\n}\n"}
